# Update the "Förändrad" date column (C2:C28) from 45526 to 45527
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45526) {
        $cell.Value = 45527
    }
}
